$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header text "MODEL_CONDITION" -> "MODELCONDITION" (currently in E1)
$ws.Range("E1").Value = "MODELCONDITION"

# Remove column A entirely (data shifts left by one column)
$ws.Range("A:A").Delete()
